$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SprintBacklog1")
$ws.Activate()
Write-Host "Windows count:" ([string]$wb.Windows.Count)
$win = $wb.Windows.Item(1)
$win.ScrollRow = 14
$win.ScrollColumn = 1
Write-Host "Win ScrollRow:" ([string]$win.ScrollRow)
$ws.Range("E20").Select()
